$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 589 (shifts existing rows 589-644 down to 590-645)
$ws.Rows("589:589").Insert()

# Populate the new row 589 with a new data record (same fixed columns as the
# surrounding "Vega Modelo de Temuco" / Mango rows, with the new record's
# date, volume, prices, origin and $/Kg values)
$ws.Range("A589").Value = 10
$ws.Range("B589").Value = "Vega Modelo de Temuco"
$ws.Range("C589").Value = "La Araucanía"
$ws.Range("D589").Value = 45194
$ws.Range("E589").Value = 9
$ws.Range("F589").Value = "Fruta"
$ws.Range("G589").Value = 100108
$ws.Range("H589").Value = "Tropicales y subtropicales"
$ws.Range("I589").Value = 100108002
$ws.Range("J589").Value = "Mango"
$ws.Range("K589").Value = "Sin especificar"
$ws.Range("L589").Value = "Primera"
$ws.Range("M589").Value = 405
$ws.Range("N589").Value = 10000
$ws.Range("O589").Value = 12000
$ws.Range("P589").Value = 11383
$ws.Range("Q589").Value = "$/bandeja 4 kilos"
$ws.Range("R589").Value = "Brasil"
$ws.Range("S589").Value = 2846
$ws.Range("T589").Value = 4
